$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows above the existing "MISCELLANEOUS" section (row 137),
# pushing it down to rows 145-148 to make room for the new
# "SOUND AND MUSIC" section documenting the v1.1 sound instructions.
$ws.Rows("137:144").Insert()

# --- Section header row (bold + underline, like the other section titles) ---
$hdr = $ws.Range("A137")
$hdr.Value = "SOUND AND MUSIC"
$hdr.NumberFormat = "@"
$hdr.Font.Bold = $true
$hdr.Font.Underline = 2

# --- New instruction rows (values written in authoring order) ---
$ws.Range("A139").Value = "b0"
$ws.Range("E142").Value = "playsnd xx yyyy"
$ws.Range("A140").Value = "b1"
$ws.Range("E143").Value = "stopsnd xx"
$ws.Range("E140").Value = "loadsnd xxxx"
$ws.Range("F140").Value = "Load sound stored in the page at address xxxx, setting up the pattern list"
$ws.Range("E139").Value = "loadsnd xx"
$ws.Range("F139").Value = "Load sound stored in the page at pointer xxxx, setting up the pattern list"
$ws.Range("A141").Value = "b2"
$ws.Range("E141").Value = "playsnd xx yy"
$ws.Range("F141").Value = "Play the pattern indexed by register yy using the playhead given by register xx (0 or 1)"
$ws.Range("A142").Value = "b3"
$ws.Range("F142").Value = "Play pattern yyyy using the playhead given by register xx (0 or 1)"
$ws.Range("A143").Value = "b4"
$ws.Range("F143").Value = "Stop the playhead given by register xx (0 or 1)"

# --- Column header row (bold, like the other column header rows) ---
$c = $ws.Range("A138")
$c.Value = "Instruction"
$c.NumberFormat = "@"
$c.Font.Bold = $true

$c = $ws.Range("E138")
$c.Value = "Assembly"
$c.NumberFormat = "@"
$c.Font.Bold = $true

$c = $ws.Range("F138")
$c.Value = "Description"
$c.NumberFormat = "@"
$c.Font.Bold = $true

# --- Remaining operand columns (reuse existing shared strings: "00", "xx", "yy") ---
$ws.Range("B139").Value = "00"
$ws.Range("C139").Value = "xx"
$ws.Range("D139").Value = "00"
$ws.Range("B140").Value = "00"
$ws.Range("C140").Value = "xx"
$ws.Range("D140").Value = "xx"
$ws.Range("B141").Value = "xx"
$ws.Range("C141").Value = "yy"
$ws.Range("D141").Value = "00"
$ws.Range("B142").Value = "xx"
$ws.Range("C142").Value = "yy"
$ws.Range("D142").Value = "yy"
$ws.Range("B143").Value = "xx"
$ws.Range("C143").Value = "00"
$ws.Range("D143").Value = "00"

# Leave the selection where the author was last working.
[void]$ws.Range("E141").Select()
